$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 13).Value = "nan"
}
